$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7460932908215167
$ws.Cells.Item(2, 3).Value = 1.495811809543383
$ws.Cells.Item(3, 2).Value = 11.82099694545185
$ws.Cells.Item(3, 3).Value = 3.160810073292703
$ws.Cells.Item(4, 2).Value = 12.67176373054251
$ws.Cells.Item(4, 3).Value = 4.829723584731859
$ws.Cells.Item(5, 2).Value = 12.95651089723617
$ws.Cells.Item(5, 3).Value = 6.449940604911695
$ws.Cells.Item(6, 2).Value = 14.23430611335822
$ws.Cells.Item(6, 3).Value = 8.226664884342259
$ws.Cells.Item(7, 2).Value = 16.22560760052146
$ws.Cells.Item(7, 3).Value = 9.969156279425475
$ws.Cells.Item(8, 2).Value = 16.78328070220875
$ws.Cells.Item(8, 3).Value = 11.51781926214407
$ws.Cells.Item(9, 2).Value = 18.35236083029099
$ws.Cells.Item(9, 3).Value = 13.02751198800834
$ws.Cells.Item(10, 2).Value = 19.29107040345616
$ws.Cells.Item(10, 3).Value = 14.70912995463671
$ws.Cells.Item(11, 2).Value = 22.01654219451287
$ws.Cells.Item(11, 3).Value = 16.30928092003757
$ws.Cells.Item(12, 2).Value = 24.3878467458843
$ws.Cells.Item(12, 3).Value = 18.11931684300703
$ws.Cells.Item(13, 2).Value = 24.88758684916857
$ws.Cells.Item(13, 3).Value = 19.78957274574999
$ws.Cells.Item(14, 2).Value = 29.27158397278049
$ws.Cells.Item(14, 3).Value = 21.37381201218944
$ws.Cells.Item(15, 2).Value = 34.39936865213679
$ws.Cells.Item(15, 3).Value = 22.98183979676146
$ws.Cells.Item(16, 2).Value = 34.94233597961408
$ws.Cells.Item(16, 3).Value = 24.69173610457452
$ws.Cells.Item(17, 2).Value = 36.32422648444486
$ws.Cells.Item(17, 3).Value = 26.29822353286865
$ws.Cells.Item(18, 2).Value = 38.13355932778981
$ws.Cells.Item(18, 3).Value = 28.0665786138402
$ws.Cells.Item(19, 2).Value = 40.76623283660203
$ws.Cells.Item(19, 3).Value = 29.66989491055532
$ws.Cells.Item(20, 2).Value = 42.0815696278426
$ws.Cells.Item(20, 3).Value = 31.47201061719143
$ws.Cells.Item(21, 2).Value = 49.22777251594825
$ws.Cells.Item(21, 3).Value = 32.93921049218886
$ws.Cells.Item(22, 2).Value = 51.00453631538158
$ws.Cells.Item(22, 3).Value = 35.14059556524872
$ws.Cells.Item(23, 2).Value = 51.32463410021728
$ws.Cells.Item(23, 3).Value = 36.81924349588125
$ws.Cells.Item(24, 2).Value = 52.78562896435605
$ws.Cells.Item(24, 3).Value = 38.73890303768326
$ws.Cells.Item(25, 2).Value = 55.70576295212567
$ws.Cells.Item(25, 3).Value = 40.55343750710154
$ws.Cells.Item(26, 2).Value = 57.34215450740146
$ws.Cells.Item(26, 3).Value = 42.35486234004176
$ws.Cells.Item(27, 2).Value = 62.24864743419364
$ws.Cells.Item(27, 3).Value = 44.03161235288869
$ws.Cells.Item(28, 2).Value = 66.52937133250907
$ws.Cells.Item(28, 3).Value = 45.60308546056995
$ws.Cells.Item(29, 2).Value = 68.1659268928539
$ws.Cells.Item(29, 3).Value = 47.36781281737079
$ws.Cells.Item(30, 2).Value = 69.45240930068093
$ws.Cells.Item(30, 3).Value = 49.17779891995099
$ws.Cells.Item(31, 2).Value = 70.0956272008829
$ws.Cells.Item(31, 3).Value = 50.79309677605689
$ws.Cells.Item(32, 2).Value = 70.40706527771319
$ws.Cells.Item(32, 3).Value = 52.2934049211719
$ws.Cells.Item(33, 2).Value = 72.33771536605643
$ws.Cells.Item(33, 3).Value = 53.79565065367358
$ws.Cells.Item(34, 2).Value = 73.38277289202919
$ws.Cells.Item(34, 3).Value = 55.41373920242348
$ws.Cells.Item(35, 2).Value = 75.50074034551564
$ws.Cells.Item(35, 3).Value = 57.06841212270529
$ws.Cells.Item(36, 2).Value = 77.3358726272329
$ws.Cells.Item(36, 3).Value = 58.56839825961777
$ws.Cells.Item(37, 2).Value = 80.18632198165018
$ws.Cells.Item(37, 3).Value = 60.53483666688842
$ws.Cells.Item(38, 2).Value = 80.73797171844051
$ws.Cells.Item(38, 3).Value = 62.32084091449345
$ws.Cells.Item(39, 2).Value = 83.47839930979748
$ws.Cells.Item(39, 3).Value = 64.07487183488219
$ws.Cells.Item(40, 2).Value = 85.73130143419417
$ws.Cells.Item(40, 3).Value = 65.59658717386358
$ws.Cells.Item(41, 2).Value = 86.37656856841075
$ws.Cells.Item(41, 3).Value = 67.1644519209344
$ws.Cells.Item(42, 2).Value = 86.7979350073706
$ws.Cells.Item(42, 3).Value = 68.67014256239486
$ws.Cells.Item(43, 2).Value = 91.31469147667551
$ws.Cells.Item(43, 3).Value = 70.57451601455645
$ws.Cells.Item(44, 2).Value = 92.79587986962537
$ws.Cells.Item(44, 3).Value = 72.79976150410505
$ws.Cells.Item(45, 2).Value = 95.45993842329972
$ws.Cells.Item(45, 3).Value = 74.54725591913616
$ws.Cells.Item(46, 2).Value = 96.18789627083585
$ws.Cells.Item(46, 3).Value = 76.45073886221356
$ws.Cells.Item(47, 2).Value = 98.62286575930176
$ws.Cells.Item(47, 3).Value = 77.95455279832527
$ws.Cells.Item(48, 2).Value = 99.56382353471527
$ws.Cells.Item(48, 3).Value = 79.80071399837608
